$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need to be marked (filtered/matched rows): 4-7 and 10-17
$rows = @(4,5,6,7,10,11,12,13,14,15,16,17)

foreach ($r in $rows) {
    # Apply the same style as the already-highlighted rows (yellow fill) across A:K
    $rng = $ws.Range("A" + $r + ":K" + $r)
    $rng.Interior.Color = 65535  # yellow (RGB 255,255,0 -> BGR 65535)

    # Set column K (conciliada) to 1 for these rows
    $ws.Cells.Item($r, 11).Value = 1
}
